$d = $word.ActiveDocument

# Smart-quote characters matching the rest of the document's prose.
$apos = [char]0x2019
$ldq  = [char]0x201C
$rdq  = [char]0x201D

# ---------------------------------------------------------------------
# 1) Find the obsolete tail of the sentence:
#    " notebook extension to easily see the charts. If this doesn't work
#    for you I will include all the same material at the bottom of this "
# and strip everything from " If this doesn't..." onward, leaving just
# " notebook extension to easily see the charts." in place.
# ---------------------------------------------------------------------
$oldTail = " If this doesn" + $apos + "t work for you I will include all the same material at the bottom of this "

$rng = $d.Content
$found = $rng.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target sentence to replace."
}

$rng.Text = ""
$rng.Collapse(0)

# ---------------------------------------------------------------------
# 2) Append the new continuation sentence right after "charts."
# ---------------------------------------------------------------------
$rng.InsertAfter(" As long as dataAnalysis.py is opened in ")
$rng.Collapse(0)

$rng.InsertAfter("VScode")
$rng.Collapse(0)

$rng.InsertAfter(" with the ")
$rng.Collapse(0)

$rng.InsertAfter("jupyter")
$rng.Collapse(0)

$closing = " notebook extension enabled, the analysis should be viewable by clicking " + $ldq + "Run Cell" + $rdq + " at the top of the file."
$rng.InsertAfter($closing)
$rng.Collapse(0)

# ---------------------------------------------------------------------
# 3) Split the paragraph right here. The remainder of the original
# paragraph ("README" + ".") moves into a brand-new paragraph that
# inherits the same "No Spacing" style.
# ---------------------------------------------------------------------
$splitPos = $rng.End
$rng.InsertParagraphAfter()

$lead = $d.Range($splitPos + 1, $splitPos + 1)

# ---------------------------------------------------------------------
# 4) Insert the new lead-in sentence in front of the leftover
# "README." runs, matching the target text.
# ---------------------------------------------------------------------
$lead.InsertAfter("Just in case that fails to work")
$lead.Collapse(0)

$lead.InsertAfter(" I ")
$lead.Collapse(0)

$lead.InsertAfter("have")
$lead.Collapse(0)

$lead.InsertAfter(" include")
$lead.Collapse(0)

$lead.InsertAfter("d")
$lead.Collapse(0)

$lead.InsertAfter(" all the same material at the bottom of this ")
$lead.Collapse(0)

Write-Output "edit applied"
